# Timesheet update: new week (09-26 to 09-30), drop placeholder tasks,
# add "Research CDS" / "Practice React", log hours, refresh totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the week's dates (row 5) ---
$ws.Range("B5").Value = 44830
$ws.Range("C5").Value = 44831
$ws.Range("D5").Value = 44832
$ws.Range("E5").Value = 42033
$ws.Range("F5").Value = 44834
$ws.Range("G5").Value = 44835
$ws.Range("H5").Value = 44836

# --- Log hours on existing task rows ---
$ws.Range("B6").Value = 1

$ws.Range("B8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = 1

# --- Remove the four placeholder task rows (Task ccccc/xxxxxx/yyyyy/zzzzzz) ---
$ws.Rows("12:15").Delete()

# --- Rename remaining task rows and log hours for the new tasks ---
$ws.Range("A10").Value = "Research CDS"
$ws.Range("B10:G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("B10").Copy()
$ws.Range("H10").PasteSpecial(-4122)

$ws.Range("A11").Value = "Practice React"
$ws.Range("B11:G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("B11").Copy()
$ws.Range("H11").PasteSpecial(-4122)

# --- Restore shared-formula grouping for the daily-total column ---
$ws.Range("I6:I11").Formula = "=SUM(B6:H6)"

# --- Column width tweaks ---
$ws.Columns("B:D").ColumnWidth = 6.166666666666667
$ws.Columns("F:F").ColumnWidth = 6.166666666666667

# --- Selection as left by the author ---
$ws.Range("L8").Select()
